$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell value while forcing Text storage (prevents Excel
# from auto-coercing numeric-looking strings like "20.10" or "0.9970"
# into Number cells, which would silently drop trailing zeros / switch
# to scientific notation for very small values such as "0.00001043").
function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

# Row 2
Set-TextValue $ws.Range("D2") '27.206.97'
Set-TextValue $ws.Range("E2") '  +0.39%  '

# Row 3
Set-TextValue $ws.Range("D3") '1.769.14'
Set-TextValue $ws.Range("E3") '  +3.04%  '

# Row 4
Set-TextValue $ws.Range("D4") '0.9971'
Set-TextValue $ws.Range("E4") '  -0.44%  '

# Row 5
Set-TextValue $ws.Range("D5") '313.08'
Set-TextValue $ws.Range("E5") '  +1.55%  '

# Row 6
Set-TextValue $ws.Range("D6") '0.9970'
Set-TextValue $ws.Range("E6") '  -0.40%  '

# Row 7
Set-TextValue $ws.Range("D7") '0.5183'
Set-TextValue $ws.Range("E7") '  +9.67%  '

# Row 8
Set-TextValue $ws.Range("D8") '0.3610'
Set-TextValue $ws.Range("E8") '  +5.09%  '

# Row 9
Set-TextValue $ws.Range("D9") '42.35'
Set-TextValue $ws.Range("E9") '  +0.57%  '

# Row 10
Set-TextValue $ws.Range("D10") '0.07329'
Set-TextValue $ws.Range("E10") '  +0.69%  '

# Row 11
Set-TextValue $ws.Range("D11") '1.083'
Set-TextValue $ws.Range("E11") '  +3.73%  '

# Row 12
Set-TextValue $ws.Range("D12") '0.9964'
Set-TextValue $ws.Range("E12") '  -0.45%  '

# Row 13
Set-TextValue $ws.Range("D13") '20.51'
Set-TextValue $ws.Range("E13") '  +2.91%  '

# Row 14
Set-TextValue $ws.Range("D14") '6.044'
Set-TextValue $ws.Range("E14") '  +2.77%  '

# Row 15
Set-TextValue $ws.Range("D15") '1.761.81'
Set-TextValue $ws.Range("E15") '  +2.52%  '

# Row 16
Set-TextValue $ws.Range("D16") '6.948'
Set-TextValue $ws.Range("E16") '  +0.80%  '

# Row 17
Set-TextValue $ws.Range("D17") '88.26'
Set-TextValue $ws.Range("E17") '  -1.05%  '

# Row 18
Set-TextValue $ws.Range("D18") '0.00001043'
Set-TextValue $ws.Range("E18") '  +0.17%  '

# Row 19
Set-TextValue $ws.Range("D19") '0.06416'
Set-TextValue $ws.Range("E19") '  +0.90%  '

# Row 20
Set-TextValue $ws.Range("D20") '0.9967'
Set-TextValue $ws.Range("E20") '  -0.41%  '

# Row 21
Set-TextValue $ws.Range("D21") '16.74'
Set-TextValue $ws.Range("E21") '  +1.24%  '

# Row 22
Set-TextValue $ws.Range("D22") '5.822'
Set-TextValue $ws.Range("E22") '  +3.39%  '

# Row 23
Set-TextValue $ws.Range("D23") '27.289.03'
Set-TextValue $ws.Range("E23") '  +0.57%  '

# Row 24
Set-TextValue $ws.Range("D24") '11.35'
Set-TextValue $ws.Range("E24") '  +4.47%  '

# Row 25
Set-TextValue $ws.Range("D25") '2.059'
Set-TextValue $ws.Range("E25") '  -3.01%  '

# Row 26
Set-TextValue $ws.Range("D26") '154.56'
Set-TextValue $ws.Range("E26") '  -1.35%  '

# Row 27
Set-TextValue $ws.Range("D27") '20.10'
Set-TextValue $ws.Range("E27") '  +2.88%  '

# Row 28
Set-TextValue $ws.Range("B28") 'LidoDAOToken'
Set-TextValue $ws.Range("C28") 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
Set-TextValue $ws.Range("D28") '2.315'
Set-TextValue $ws.Range("E28") '  +10.41%  '

# Row 29
Set-TextValue $ws.Range("B29") 'WrappedliquidstakedEther2.0'
Set-TextValue $ws.Range("C29") 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
Set-TextValue $ws.Range("D29") '1.956.89'
Set-TextValue $ws.Range("E29") '  +2.53%  '

# Row 30
Set-TextValue $ws.Range("D30") '121.14'
Set-TextValue $ws.Range("E30") '  +1.26%  '

# Row 31
Set-TextValue $ws.Range("D31") '1.059'
Set-TextValue $ws.Range("E31") '  +4.06%  '

# Row 32
Set-TextValue $ws.Range("D32") '0.09689'
Set-TextValue $ws.Range("E32") '  +5.68%  '

# Row 33
Set-TextValue $ws.Range("B33") 'Filecoin'
Set-TextValue $ws.Range("C33") 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-TextValue $ws.Range("D33") '5.506'
Set-TextValue $ws.Range("E33") '  +3.28%  '

# Row 34
Set-TextValue $ws.Range("B34") 'HuobiToken'
Set-TextValue $ws.Range("C34") 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
Set-TextValue $ws.Range("D34") '3.585'
Set-TextValue $ws.Range("E34") '  -0.23%  '

# Row 35
Set-TextValue $ws.Range("B35") 'VeChain'
Set-TextValue $ws.Range("C35") 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextValue $ws.Range("D35") '0.02219'
Set-TextValue $ws.Range("E35") '  +0.53%  '

# Row 36
Set-TextValue $ws.Range("B36") 'Hedera'
Set-TextValue $ws.Range("C36") 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-TextValue $ws.Range("D36") '0.05985'
Set-TextValue $ws.Range("E36") '  +2.59%  '

# Row 37
Set-TextValue $ws.Range("D37") '11.20'
Set-TextValue $ws.Range("E37") '  +1.71%  '

# Row 38
Set-TextValue $ws.Range("D38") '0.2027'
Set-TextValue $ws.Range("E38") '  +1.37%  '

# Row 39
Set-TextValue $ws.Range("D39") '4.826'
Set-TextValue $ws.Range("E39") '  +1.70%  '

# Row 40
Set-TextValue $ws.Range("D40") '0.6115'
Set-TextValue $ws.Range("E40") '  +3.65%  '

# Row 41
Set-TextValue $ws.Range("D41") '1.431'
Set-TextValue $ws.Range("E41") '  +2.73%  '

# Row 42
Set-TextValue $ws.Range("B42") 'TrustWalletToken'
Set-TextValue $ws.Range("C42") 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
Set-TextValue $ws.Range("D42") '1.125'
Set-TextValue $ws.Range("E42") '  +0.23%  '

# Row 43
Set-TextValue $ws.Range("B43") 'FraxShare'
Set-TextValue $ws.Range("C43") 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
Set-TextValue $ws.Range("D43") '7.818'
Set-TextValue $ws.Range("E43") '  +4.65%  '

# Row 44
Set-TextValue $ws.Range("D44") '13.12'
Set-TextValue $ws.Range("E44") '  +4.19%  '

# Row 45
Set-TextValue $ws.Range("D45") '3.620'
Set-TextValue $ws.Range("E45") '  +1.66%  '

# Row 46
Set-TextValue $ws.Range("D46") '0.5735'
Set-TextValue $ws.Range("E46") '  +1.24%  '

# Row 47
Set-TextValue $ws.Range("D47") '121.43'
Set-TextValue $ws.Range("E47") '  +3.19%  '

# Row 48
Set-TextValue $ws.Range("D48") '1.877'
Set-TextValue $ws.Range("E48") '  +1.83%  '

# Row 49
Set-TextValue $ws.Range("D49") '1.110'
Set-TextValue $ws.Range("E49") '  +2.18%  '

# Row 50
Set-TextValue $ws.Range("D50") '0.06697'
Set-TextValue $ws.Range("E50") '  +0.72%  '

# Row 51
Set-TextValue $ws.Range("D51") '70.47'
Set-TextValue $ws.Range("E51") '  +0.72%  '
